$d = $word.ActiveDocument

$d.Content.Find.Execute("4+16=", $true, $false, $false, $false, $false, $true, 1, $false, "57-48=", 2) | Out-Null
$d.Content.Find.Execute("88-81=", $true, $false, $false, $false, $false, $true, 1, $false, "66-27=", 2) | Out-Null
$d.Content.Find.Execute("86+0=", $true, $false, $false, $false, $false, $true, 1, $false, "9+67=", 2) | Out-Null
$d.Content.Find.Execute("21+20=", $true, $false, $false, $false, $false, $true, 1, $false, "81-69=", 2) | Out-Null
$d.Content.Find.Execute("34+18=", $true, $false, $false, $false, $false, $true, 1, $false, "30+28=", 2) | Out-Null
$d.Content.Find.Execute("50-33=", $true, $false, $false, $false, $false, $true, 1, $false, "10+43=", 2) | Out-Null
$d.Content.Find.Execute("6+67=", $true, $false, $false, $false, $false, $true, 1, $false, "51+28=", 2) | Out-Null
$d.Content.Find.Execute("6+92=", $true, $false, $false, $false, $false, $true, 1, $false, "30+40=", 2) | Out-Null
$d.Content.Find.Execute("86+3=", $true, $false, $false, $false, $false, $true, 1, $false, "47-12=", 2) | Out-Null
$d.Content.Find.Execute("53-46=", $true, $false, $false, $false, $false, $true, 1, $false, "40-7=", 2) | Out-Null
$d.Content.Find.Execute("24-12=", $true, $false, $false, $false, $false, $true, 1, $false, "63+17=", 2) | Out-Null
$d.Content.Find.Execute("78-29=", $true, $false, $false, $false, $false, $true, 1, $false, "41-12=", 2) | Out-Null
$d.Content.Find.Execute("83-15=", $true, $false, $false, $false, $false, $true, 1, $false, "70-0=", 2) | Out-Null
$d.Content.Find.Execute("67-54=", $true, $false, $false, $false, $false, $true, 1, $false, "14+61=", 2) | Out-Null
$d.Content.Find.Execute("44+33=", $true, $false, $false, $false, $false, $true, 1, $false, "32-26=", 2) | Out-Null
$d.Content.Find.Execute("90-84=", $true, $false, $false, $false, $false, $true, 1, $false, "21-9=", 2) | Out-Null
$d.Content.Find.Execute("49+34=", $true, $false, $false, $false, $false, $true, 1, $false, "58-58=", 2) | Out-Null
$d.Content.Find.Execute("34-16=", $true, $false, $false, $false, $false, $true, 1, $false, "44-5=", 2) | Out-Null
$d.Content.Find.Execute("31-28=", $true, $false, $false, $false, $false, $true, 1, $false, "12+21=", 2) | Out-Null
$d.Content.Find.Execute("63-21=", $true, $false, $false, $false, $false, $true, 1, $false, "64+3=", 2) | Out-Null
$d.Content.Find.Execute("83-36=", $true, $false, $false, $false, $false, $true, 1, $false, "80-66=", 2) | Out-Null
$d.Content.Find.Execute("25-8=", $true, $false, $false, $false, $false, $true, 1, $false, "74+17=", 2) | Out-Null
$d.Content.Find.Execute("39+10=", $true, $false, $false, $false, $false, $true, 1, $false, "6+64=", 2) | Out-Null
$d.Content.Find.Execute("2+13=", $true, $false, $false, $false, $false, $true, 1, $false, "65+29=", 2) | Out-Null
$d.Content.Find.Execute("57-9=", $true, $false, $false, $false, $false, $true, 1, $false, "78+21=", 2) | Out-Null
$d.Content.Find.Execute("16+80=", $true, $false, $false, $false, $false, $true, 1, $false, "87-49=", 2) | Out-Null
$d.Content.Find.Execute("83-69=", $true, $false, $false, $false, $false, $true, 1, $false, "12+20=", 2) | Out-Null
$d.Content.Find.Execute("95-82=", $true, $false, $false, $false, $false, $true, 1, $false, "97-51=", 2) | Out-Null
$d.Content.Find.Execute("60+35=", $true, $false, $false, $false, $false, $true, 1, $false, "31-30=", 2) | Out-Null
$d.Content.Find.Execute("0+4=", $true, $false, $false, $false, $false, $true, 1, $false, "49+49=", 2) | Out-Null
$d.Content.Find.Execute("82+0=", $true, $false, $false, $false, $false, $true, 1, $false, "63-41=", 2) | Out-Null
$d.Content.Find.Execute("19-14=", $true, $false, $false, $false, $false, $true, 1, $false, "2+1=", 2) | Out-Null
$d.Content.Find.Execute("12+73=", $true, $false, $false, $false, $false, $true, 1, $false, "63-50=", 2) | Out-Null
$d.Content.Find.Execute("41+3=", $true, $false, $false, $false, $false, $true, 1, $false, "82-59=", 2) | Out-Null
$d.Content.Find.Execute("38+47=", $true, $false, $false, $false, $false, $true, 1, $false, "35-2=", 2) | Out-Null
$d.Content.Find.Execute("37-29=", $true, $false, $false, $false, $false, $true, 1, $false, "61-50=", 2) | Out-Null
$d.Content.Find.Execute("85-16=", $true, $false, $false, $false, $false, $true, 1, $false, "76-39=", 2) | Out-Null
$d.Content.Find.Execute("52+40=", $true, $false, $false, $false, $false, $true, 1, $false, "1+12=", 2) | Out-Null
$d.Content.Find.Execute("38-26=", $true, $false, $false, $false, $false, $true, 1, $false, "8+77=", 2) | Out-Null
$d.Content.Find.Execute("42-8=", $true, $false, $false, $false, $false, $true, 1, $false, "0+70=", 2) | Out-Null
$d.Content.Find.Execute("51-44=", $true, $false, $false, $false, $false, $true, 1, $false, "21+55=", 2) | Out-Null
$d.Content.Find.Execute("60-13=", $true, $false, $false, $false, $false, $true, 1, $false, "3+12=", 2) | Out-Null
$d.Content.Find.Execute("66-59=", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=", 2) | Out-Null
$d.Content.Find.Execute("4+45=", $true, $false, $false, $false, $false, $true, 1, $false, "5+56=", 2) | Out-Null
$d.Content.Find.Execute("61+1=", $true, $false, $false, $false, $false, $true, 1, $false, "72-54=", 2) | Out-Null
$d.Content.Find.Execute("26-25=", $true, $false, $false, $false, $false, $true, 1, $false, "67+19=", 2) | Out-Null
$d.Content.Find.Execute("4+40=", $true, $false, $false, $false, $false, $true, 1, $false, "38+38=", 2) | Out-Null
$d.Content.Find.Execute("23+46=", $true, $false, $false, $false, $false, $true, 1, $false, "34+29=", 2) | Out-Null
$d.Content.Find.Execute("52+0=", $true, $false, $false, $false, $false, $true, 1, $false, "57-49=", 2) | Out-Null
$d.Content.Find.Execute("68-11=", $true, $false, $false, $false, $false, $true, 1, $false, "10+41=", 2) | Out-Null
$d.Content.Find.Execute("51-3=", $true, $false, $false, $false, $false, $true, 1, $false, "67+28=", 2) | Out-Null
$d.Content.Find.Execute("10+0=", $true, $false, $false, $false, $false, $true, 1, $false, "14+3=", 2) | Out-Null
$d.Content.Find.Execute("58-17=", $true, $false, $false, $false, $false, $true, 1, $false, "14+55=", 2) | Out-Null
$d.Content.Find.Execute("90-49=", $true, $false, $false, $false, $false, $true, 1, $false, "26+64=", 2) | Out-Null
$d.Content.Find.Execute("81-17=", $true, $false, $false, $false, $false, $true, 1, $false, "59-24=", 2) | Out-Null
$d.Content.Find.Execute("22+4=", $true, $false, $false, $false, $false, $true, 1, $false, "81-60=", 2) | Out-Null
$d.Content.Find.Execute("69+9=", $true, $false, $false, $false, $false, $true, 1, $false, "32+56=", 2) | Out-Null
$d.Content.Find.Execute("74-71=", $true, $false, $false, $false, $false, $true, 1, $false, "39-25=", 2) | Out-Null
$d.Content.Find.Execute("24+52=", $true, $false, $false, $false, $false, $true, 1, $false, "75+22=", 2) | Out-Null
$d.Content.Find.Execute("63+24=", $true, $false, $false, $false, $false, $true, 1, $false, "33+32=", 2) | Out-Null
$d.Content.Find.Execute("9+29=", $true, $false, $false, $false, $false, $true, 1, $false, "25+18=", 2) | Out-Null
$d.Content.Find.Execute("30+9=", $true, $false, $false, $false, $false, $true, 1, $false, "72-43=", 2) | Out-Null
$d.Content.Find.Execute("35+61=", $true, $false, $false, $false, $false, $true, 1, $false, "79-75=", 2) | Out-Null
$d.Content.Find.Execute("35-9=", $true, $false, $false, $false, $false, $true, 1, $false, "94-11=", 2) | Out-Null
$d.Content.Find.Execute("46+43=", $true, $false, $false, $false, $false, $true, 1, $false, "87-65=", 2) | Out-Null
$d.Content.Find.Execute("82-53=", $true, $false, $false, $false, $false, $true, 1, $false, "96-86=", 2) | Out-Null
$d.Content.Find.Execute("96-21=", $true, $false, $false, $false, $false, $true, 1, $false, "35+2=", 2) | Out-Null
$d.Content.Find.Execute("38+49=", $true, $false, $false, $false, $false, $true, 1, $false, "36+14=", 2) | Out-Null
$d.Content.Find.Execute("20+39=", $true, $false, $false, $false, $false, $true, 1, $false, "37-18=", 2) | Out-Null
$d.Content.Find.Execute("10+66=", $true, $false, $false, $false, $false, $true, 1, $false, "63-33=", 2) | Out-Null
$d.Content.Find.Execute("95-67=", $true, $false, $false, $false, $false, $true, 1, $false, "29+15=", 2) | Out-Null
$d.Content.Find.Execute("20-7=", $true, $false, $false, $false, $false, $true, 1, $false, "32+27=", 2) | Out-Null
$d.Content.Find.Execute("74-5=", $true, $false, $false, $false, $false, $true, 1, $false, "46+12=", 2) | Out-Null
$d.Content.Find.Execute("76-34=", $true, $false, $false, $false, $false, $true, 1, $false, "24+7=", 2) | Out-Null
$d.Content.Find.Execute("15-0=", $true, $false, $false, $false, $false, $true, 1, $false, "26-15=", 2) | Out-Null
$d.Content.Find.Execute("82-71=", $true, $false, $false, $false, $false, $true, 1, $false, "10+79=", 2) | Out-Null
$d.Content.Find.Execute("92-37=", $true, $false, $false, $false, $false, $true, 1, $false, "64-30=", 2) | Out-Null
$d.Content.Find.Execute("30-12=", $true, $false, $false, $false, $false, $true, 1, $false, "66-6=", 2) | Out-Null
$d.Content.Find.Execute("4+41=", $true, $false, $false, $false, $false, $true, 1, $false, "18+32=", 2) | Out-Null
$d.Content.Find.Execute("27-9=", $true, $false, $false, $false, $false, $true, 1, $false, "2+95=", 2) | Out-Null
$d.Content.Find.Execute("87+3=", $true, $false, $false, $false, $false, $true, 1, $false, "19-15=", 2) | Out-Null
$d.Content.Find.Execute("46+25=", $true, $false, $false, $false, $false, $true, 1, $false, "44-3=", 2) | Out-Null
$d.Content.Find.Execute("14+8=", $true, $false, $false, $false, $false, $true, 1, $false, "81+5=", 2) | Out-Null
$d.Content.Find.Execute("94-10=", $true, $false, $false, $false, $false, $true, 1, $false, "38+42=", 2) | Out-Null
$d.Content.Find.Execute("73-62=", $true, $false, $false, $false, $false, $true, 1, $false, "65-21=", 2) | Out-Null
$d.Content.Find.Execute("56-42=", $true, $false, $false, $false, $false, $true, 1, $false, "34-30=", 2) | Out-Null
$d.Content.Find.Execute("13+19=", $true, $false, $false, $false, $false, $true, 1, $false, "90-71=", 2) | Out-Null
$d.Content.Find.Execute("69-41=", $true, $false, $false, $false, $false, $true, 1, $false, "38+32=", 2) | Out-Null
$d.Content.Find.Execute("65+12=", $true, $false, $false, $false, $false, $true, 1, $false, "62-38=", 2) | Out-Null
$d.Content.Find.Execute("4+79=", $true, $false, $false, $false, $false, $true, 1, $false, "82-20=", 2) | Out-Null
$d.Content.Find.Execute("42-34=", $true, $false, $false, $false, $false, $true, 1, $false, "89-14=", 2) | Out-Null
$d.Content.Find.Execute("77-3=", $true, $false, $false, $false, $false, $true, 1, $false, "45-22=", 2) | Out-Null
$d.Content.Find.Execute("26+22=", $true, $false, $false, $false, $false, $true, 1, $false, "51-9=", 2) | Out-Null
$d.Content.Find.Execute("92-74=", $true, $false, $false, $false, $false, $true, 1, $false, "29+20=", 2) | Out-Null
$d.Content.Find.Execute("84-19=", $true, $false, $false, $false, $false, $true, 1, $false, "24+45=", 2) | Out-Null
$d.Content.Find.Execute("55-33=", $true, $false, $false, $false, $false, $true, 1, $false, "29-22=", 2) | Out-Null
$d.Content.Find.Execute("26-18=", $true, $false, $false, $false, $false, $true, 1, $false, "24+21=", 2) | Out-Null
$d.Content.Find.Execute("67-18=", $true, $false, $false, $false, $false, $true, 1, $false, "62+33=", 2) | Out-Null
$d.Content.Find.Execute("20+61=", $true, $false, $false, $false, $false, $true, 1, $false, "26+16=", 2) | Out-Null
$d.Content.Find.Execute("86-54=", $true, $false, $false, $false, $false, $true, 1, $false, "79-25=", 2) | Out-Null
